$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "94-77="
$t.Cell(1, 2).Range.Text = "9+56="
$t.Cell(1, 3).Range.Text = "19+68="
$t.Cell(1, 4).Range.Text = "54-5="
$t.Cell(1, 5).Range.Text = "63-37="
$t.Cell(2, 1).Range.Text = "31-28="
$t.Cell(2, 2).Range.Text = "39+33="
$t.Cell(2, 3).Range.Text = "45-19="
$t.Cell(2, 4).Range.Text = "70-55="
$t.Cell(2, 5).Range.Text = "87+8="
$t.Cell(3, 1).Range.Text = "27+24="
$t.Cell(3, 2).Range.Text = "29+36="
$t.Cell(3, 3).Range.Text = "52-49="
$t.Cell(3, 4).Range.Text = "50-33="
$t.Cell(3, 5).Range.Text = "10-7="
$t.Cell(4, 1).Range.Text = "28+49="
$t.Cell(4, 2).Range.Text = "9+29="
$t.Cell(4, 3).Range.Text = "49+43="
$t.Cell(4, 4).Range.Text = "52-24="
$t.Cell(4, 5).Range.Text = "90-9="
$t.Cell(5, 1).Range.Text = "19+69="
$t.Cell(5, 2).Range.Text = "9+47="
$t.Cell(5, 3).Range.Text = "80-75="
$t.Cell(5, 4).Range.Text = "48+47="
$t.Cell(5, 5).Range.Text = "65-36="
$t.Cell(6, 1).Range.Text = "37+6="
$t.Cell(6, 2).Range.Text = "50-45="
$t.Cell(6, 3).Range.Text = "60-18="
$t.Cell(6, 4).Range.Text = "82-39="
$t.Cell(6, 5).Range.Text = "87-58="
$t.Cell(7, 1).Range.Text = "50-14="
$t.Cell(7, 2).Range.Text = "42-15="
$t.Cell(7, 3).Range.Text = "90-81="
$t.Cell(7, 4).Range.Text = "27+28="
$t.Cell(7, 5).Range.Text = "41-6="
$t.Cell(8, 1).Range.Text = "84-27="
$t.Cell(8, 2).Range.Text = "70-27="
$t.Cell(8, 3).Range.Text = "62-35="
$t.Cell(8, 4).Range.Text = "67+6="
$t.Cell(8, 5).Range.Text = "84-8="
$t.Cell(9, 1).Range.Text = "23-5="
$t.Cell(9, 2).Range.Text = "59+6="
$t.Cell(9, 3).Range.Text = "81-13="
$t.Cell(9, 4).Range.Text = "34+7="
$t.Cell(9, 5).Range.Text = "58+34="
$t.Cell(10, 1).Range.Text = "90-52="
$t.Cell(10, 2).Range.Text = "14+39="
$t.Cell(10, 3).Range.Text = "32-18="
$t.Cell(10, 4).Range.Text = "28+3="
$t.Cell(10, 5).Range.Text = "30-11="
$t.Cell(11, 1).Range.Text = "61-52="
$t.Cell(11, 2).Range.Text = "83-36="
$t.Cell(11, 3).Range.Text = "28+37="
$t.Cell(11, 4).Range.Text = "94-75="
$t.Cell(11, 5).Range.Text = "71-59="
$t.Cell(12, 1).Range.Text = "7+77="
$t.Cell(12, 2).Range.Text = "74-55="
$t.Cell(12, 3).Range.Text = "90-2="
$t.Cell(12, 4).Range.Text = "96-28="
$t.Cell(12, 5).Range.Text = "49+19="
$t.Cell(13, 1).Range.Text = "57+38="
$t.Cell(13, 2).Range.Text = "35+49="
$t.Cell(13, 3).Range.Text = "56+26="
$t.Cell(13, 4).Range.Text = "16+17="
$t.Cell(13, 5).Range.Text = "4+79="
$t.Cell(14, 1).Range.Text = "37-8="
$t.Cell(14, 2).Range.Text = "3+88="
$t.Cell(14, 3).Range.Text = "72-48="
$t.Cell(14, 4).Range.Text = "94-19="
$t.Cell(14, 5).Range.Text = "19+53="
$t.Cell(15, 1).Range.Text = "4+88="
$t.Cell(15, 2).Range.Text = "40-18="
$t.Cell(15, 3).Range.Text = "80-2="
$t.Cell(15, 4).Range.Text = "26+35="
$t.Cell(15, 5).Range.Text = "16+67="
$t.Cell(16, 1).Range.Text = "42+19="
$t.Cell(16, 2).Range.Text = "90-83="
$t.Cell(16, 3).Range.Text = "26+6="
$t.Cell(16, 4).Range.Text = "7+9="
$t.Cell(16, 5).Range.Text = "18+58="
$t.Cell(17, 1).Range.Text = "70-66="
$t.Cell(17, 2).Range.Text = "93-67="
$t.Cell(17, 3).Range.Text = "5+77="
$t.Cell(17, 4).Range.Text = "53-34="
$t.Cell(17, 5).Range.Text = "58+38="
$t.Cell(18, 1).Range.Text = "9+57="
$t.Cell(18, 2).Range.Text = "91-38="
$t.Cell(18, 3).Range.Text = "6+77="
$t.Cell(18, 4).Range.Text = "25+46="
$t.Cell(18, 5).Range.Text = "28+29="
$t.Cell(19, 1).Range.Text = "28+38="
$t.Cell(19, 2).Range.Text = "52+9="
$t.Cell(19, 3).Range.Text = "92-45="
$t.Cell(19, 4).Range.Text = "86-48="
$t.Cell(19, 5).Range.Text = "15+78="
$t.Cell(20, 1).Range.Text = "26-9="
$t.Cell(20, 2).Range.Text = "9+42="
$t.Cell(20, 3).Range.Text = "83-29="
$t.Cell(20, 4).Range.Text = "66-19="
$t.Cell(20, 5).Range.Text = "79+9="
